$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 210:211, pushing the existing rows 210.. down by
# two (Excel carries formatting such as the date style on column D along).
$ws.Rows("210:211").Insert()

# New record 1 (lands on row 210 after the insert)
$row210 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44960, 9, 100112052, "Albahaca", "Sin especificar", "Primera", 30, 5000, 5000, 5000, "$/paquete", "Región de La Araucanía", 5000, 1, "Hortaliza")

# New record 2 (lands on row 211 after the insert)
$row211 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44960, 9, 100112052, "Albahaca", "Sin especificar", "Primera", 100, 4000, 4000, 4000, "$/paquete", "Región del Maule", 4000, 1, "Hortaliza")

for ($i = 0; $i -lt $row210.Length; $i++) {
    $ws.Cells.Item(210, $i + 1).Value = $row210[$i]
}

for ($i = 0; $i -lt $row211.Length; $i++) {
    $ws.Cells.Item(211, $i + 1).Value = $row211[$i]
}
